# Update the regression-output table (columns B:D, rows 2:4) with the new
# coefficient estimates (Crisis and Credit Allocation results).
#
# Some of the new values (e.g. "0.17", "-0.01") are strings that happen to
# look like plain numbers. Assigning them straight to Range.Value would make
# Excel auto-convert them into numeric cells, which is not what the source
# workbook does (every value in this table - including plain-looking ones -
# is stored as literal text). To keep them as text *without* picking up a
# quote-prefix/number-format style (which would modify styles.xml), we route
# them through a small helper cell: put a formula that evaluates to the
# desired text, copy it, and paste-special just the resulting value into the
# target cell. That yields a plain text cell identical in shape to the
# original ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values that are unambiguous text (contain non-numeric characters such as
# "*") can be assigned directly - Excel keeps them as text automatically.
$ws.Range("C2").Value = "44.29***"
$ws.Range("C3").Value = "2.21***"
$ws.Range("D3").Value = "0.46***"
$ws.Range("D4").Value = "0.82*"

# Values that look like plain numbers need the helper-cell trick so they stay
# text (matching the rest of the table) instead of becoming numeric cells.
$ws.Range("F1").Formula = '="0.17"'
$ws.Range("F1").Copy()
$ws.Range("B2").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("F1").Formula = '="-0.01"'
$ws.Range("F1").Copy()
$ws.Range("B3").PasteSpecial(-4163)

$ws.Range("F1").Formula = '="-0.09"'
$ws.Range("F1").Copy()
$ws.Range("B4").PasteSpecial(-4163)

$ws.Range("F1").Formula = '="0.98"'
$ws.Range("F1").Copy()
$ws.Range("C4").PasteSpecial(-4163)

$ws.Range("F1").Formula = '="-0.89"'
$ws.Range("F1").Copy()
$ws.Range("D2").PasteSpecial(-4163)

# Clean up the scratch cell used for the text-conversion trick.
$ws.Range("F1").ClearContents()
